# Add a new "Source" slide at the end of the deck, using the same
# Title + Content layout ("obj", slideLayout2) already used by several
# other slides, then fill in the title and body text.

$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Source"

# Body / content placeholder
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Image:"
[void]$body.InsertAfter("`rhttp://38ccda.medialib.glogster.com/media/ac57b46d05ee61023f6c297d44aaea9349af2f721b70e868bf78fda40a49fbb3/cartoon-bridge.jpg")

# Indent the URL line (2nd paragraph) one level, like the source deck.
$body.Paragraphs(2).IndentLevel = 2
